$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 203
$ws.Range("F3").Value = 9
$ws.Range("F5").Value = 177
$ws.Range("F6").Value = 534
$ws.Range("F7").Value = 41
$ws.Range("F8").Value = 9632
$ws.Range("F10").Value = 2593
$ws.Range("F11").Value = 202
$ws.Range("F12").Value = 2375
$ws.Range("F13").Value = 2610
$ws.Range("F15").Value = 267
$ws.Range("F16").Value = 2041
$ws.Range("F18").Value = 71
$ws.Range("F19").Value = 359
$ws.Range("F21").Value = 59
$ws.Range("F22").Value = 292
$ws.Range("F23").Value = 55
$ws.Range("F24").Value = 122
$ws.Range("F26").Value = 1266
$ws.Range("F27").Value = 1236
$ws.Range("F28").Value = 88
$ws.Range("F29").Value = 115
$ws.Range("F31").Value = 1635
$ws.Range("F32").Value = 2724
$ws.Range("F34").Value = 966
$ws.Range("F35").Value = 337
$ws.Range("F36").Value = 1273
$ws.Range("F37").Value = 30
$ws.Range("F38").Value = 43
$ws.Range("F40").Value = 22

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 161
$ws.Range("F8").Value = 14
$ws.Range("F16").Value = 145

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 937
$ws.Range("F4").Value = 109

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 203
$ws.Range("F4").Value = 937
$ws.Range("F5").Value = 109
$ws.Range("F6").Value = 9
$ws.Range("F9").Value = 177
$ws.Range("F10").Value = 534
$ws.Range("F11").Value = 41
$ws.Range("F12").Value = 9633
$ws.Range("F13").Value = 161
$ws.Range("F15").Value = 2593
$ws.Range("F16").Value = 202
$ws.Range("F17").Value = 2375
$ws.Range("F18").Value = 2610
$ws.Range("F19").Value = 14
$ws.Range("F20").Value = 267
$ws.Range("F21").Value = 2041
$ws.Range("F23").Value = 71
$ws.Range("F24").Value = 359
$ws.Range("F26").Value = 59
$ws.Range("F27").Value = 292
$ws.Range("F28").Value = 55
$ws.Range("F29").Value = 122
$ws.Range("F31").Value = 1266
$ws.Range("F32").Value = 1236
$ws.Range("F33").Value = 88
$ws.Range("F34").Value = 115
$ws.Range("F36").Value = 1635
$ws.Range("F38").Value = 2724
$ws.Range("F39").Value = 966
$ws.Range("F41").Value = 337
$ws.Range("F45").Value = 1273
$ws.Range("F46").Value = 43
$ws.Range("F48").Value = 22
$ws.Range("F49").Value = 145
$ws.Range("F50").Value = 145

